$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Fix H34: style change from s=10 to s=6 (value stays 80)
# -----------------------------------------------------------------
$ws.Range("H4").Copy($ws.Range("H34"))
$ws.Range("H34").Value = 80

# -----------------------------------------------------------------
# Helper pattern used below for every brand-new cell: a cell that is
# being created beyond the sheet's previous used range needs to be
# "touched" with a plain value first so the engine registers it in
# the dependency graph; only then does copying formatting from a
# template cell (to pick up the right style index) and assigning the
# real value keep dependent SUM()/formula cells correctly
# recalculated.
# -----------------------------------------------------------------

# -----------------------------------------------------------------
# Row 36: new cells A36:C36
# -----------------------------------------------------------------
$ws.Range("A36").Value = 0
$ws.Range("A4").Copy($ws.Range("A36"))
$ws.Range("A36").Value = 43071

$ws.Range("B36").Value = 0
$ws.Range("B4").Copy($ws.Range("B36"))
$ws.Range("B36").Value = "Images, Debugging, and Comments"

$ws.Range("C36").Value = 0
$ws.Range("C4").Copy($ws.Range("C36"))
$ws.Range("C36").Value = 90

# -----------------------------------------------------------------
# Row 37: new cells A37:C37
# -----------------------------------------------------------------
$ws.Range("A37").Value = 0
$ws.Range("A3").Copy($ws.Range("A37"))
$ws.Range("A37").Value = 43072

$ws.Range("B37").Value = 0
$ws.Range("B3").Copy($ws.Range("B37"))
$ws.Range("B37").Value = "Debugging,Report and Comments"

$ws.Range("C37").Value = 0
$ws.Range("C3").Copy($ws.Range("C37"))
$ws.Range("C37").Value = 60

# -----------------------------------------------------------------
# Row 38: new cells A38:C38
# -----------------------------------------------------------------
$ws.Range("A38").Value = 0
$ws.Range("A4").Copy($ws.Range("A38"))
$ws.Range("A38").Value = 43073

$ws.Range("B38").Value = 0
$ws.Range("B4").Copy($ws.Range("B38"))
$ws.Range("B38").Value = "Comments and ReadMe"

$ws.Range("C38").Value = 0
$ws.Range("C4").Copy($ws.Range("C38"))
$ws.Range("C38").Value = 90

# -----------------------------------------------------------------
# Row 39: new cells A39:C39 and K39:N39
# -----------------------------------------------------------------
$ws.Range("A39").Value = 0
$ws.Range("A3").Copy($ws.Range("A39"))
$ws.Range("A39").Value = 43073

$ws.Range("B39").Value = 0
$ws.Range("B3").Copy($ws.Range("B39"))
$ws.Range("B39").Value = "Team Meeting"

$ws.Range("C39").Value = 0
$ws.Range("C35").Copy($ws.Range("C39"))
$ws.Range("C39").Value = 10

$ws.Range("K39").Value = 0
$ws.Range("K3").Copy($ws.Range("K39"))
$ws.Range("K39").Value = 43069

$ws.Range("L39").Value = 0
$ws.Range("L3").Copy($ws.Range("L39"))
$ws.Range("L39").Value = "Bug fixes, repaired unit tests, added Boss unit tests"

$ws.Range("M39").Value = 0
$ws.Range("M3").Copy($ws.Range("M39"))
$ws.Range("M39").Value = 193

$ws.Range("N39").Value = 0
$ws.Range("N39").Formula = "=60+60+60+13"

# -----------------------------------------------------------------
# Row 40: new cells K40:M40
# -----------------------------------------------------------------
$ws.Range("K40").Value = 0
$ws.Range("K4").Copy($ws.Range("K40"))
$ws.Range("K40").Value = 43071

$ws.Range("L40").Value = 0
$ws.Range("L4").Copy($ws.Range("L40"))
$ws.Range("L40").Value = "fixed loading of Speed and invincibility, and loading medium boss state"

$ws.Range("M40").Value = 0
$ws.Range("M4").Copy($ws.Range("M40"))
$ws.Range("M40").Value = 83

# -----------------------------------------------------------------
# Row 41: new cells K41:M41
# -----------------------------------------------------------------
$ws.Range("K41").Value = 0
$ws.Range("K3").Copy($ws.Range("K41"))
$ws.Range("K41").Value = 43072

$ws.Range("L41").Value = 0
$ws.Range("L3").Copy($ws.Range("L41"))
$ws.Range("L41").Value = "fixed more bugs (hard boss and loading), had brother test for bugs (got mostly complaints), added more unit tests"

$ws.Range("M41").Value = 0
$ws.Range("M3").Copy($ws.Range("M41"))
$ws.Range("M41").Value = 217

# -----------------------------------------------------------------
# Row 42: new cells K42:M42
# -----------------------------------------------------------------
$ws.Range("K42").Value = 0
$ws.Range("K4").Copy($ws.Range("K42"))
$ws.Range("K42").Value = 43072

$ws.Range("L42").Value = 0
$ws.Range("L4").Copy($ws.Range("L42"))
$ws.Range("L42").Value = "finished up unit tests"

$ws.Range("M42").Value = 0
$ws.Range("M4").Copy($ws.Range("M42"))
$ws.Range("M42").Value = 15

# -----------------------------------------------------------------
# Row 43: new cells K43:M43
# -----------------------------------------------------------------
$ws.Range("K43").Value = 0
$ws.Range("K3").Copy($ws.Range("K43"))
$ws.Range("K43").Value = 43072

$ws.Range("L43").Value = 0
$ws.Range("L3").Copy($ws.Range("L43"))
$ws.Range("L43").Value = "Wrote comments and labels"

$ws.Range("M43").Value = 0
$ws.Range("M3").Copy($ws.Range("M43"))
$ws.Range("M43").Value = 60

# -----------------------------------------------------------------
# Row 44: new cells K44:M44
# -----------------------------------------------------------------
$ws.Range("K44").Value = 0
$ws.Range("K4").Copy($ws.Range("K44"))
$ws.Range("K44").Value = 43072

$ws.Range("L44").Value = 0
$ws.Range("L4").Copy($ws.Range("L44"))
$ws.Range("L44").Value = "Headers and comments, made release video"

$ws.Range("M44").Value = 0
$ws.Range("M4").Copy($ws.Range("M44"))
$ws.Range("M44").Value = 23

# -----------------------------------------------------------------
# Row 45: new cells K45:M45
# -----------------------------------------------------------------
$ws.Range("K45").Value = 0
$ws.Range("K3").Copy($ws.Range("K45"))
$ws.Range("K45").Value = 43073

$ws.Range("L45").Value = 0
$ws.Range("L3").Copy($ws.Range("L45"))
$ws.Range("L45").Value = "Updated serialization design, personal report"

$ws.Range("M45").Value = 0
$ws.Range("M3").Copy($ws.Range("M45"))
$ws.Range("M45").Value = 35

# -----------------------------------------------------------------
# New rows 39-45 need the same explicit row height as the rest of
# the sheet (defaultRowHeight doesn't automatically get stamped on
# freshly-created rows).
# -----------------------------------------------------------------
$ws.Rows.Item(39).RowHeight = 27.55
$ws.Rows.Item(40).RowHeight = 27.55
$ws.Rows.Item(41).RowHeight = 27.55
$ws.Rows.Item(42).RowHeight = 27.55
$ws.Rows.Item(43).RowHeight = 27.55
$ws.Rows.Item(44).RowHeight = 27.55
$ws.Rows.Item(45).RowHeight = 27.55

# -----------------------------------------------------------------
# Update the view: top-left cell and selection
# -----------------------------------------------------------------
$ws.Activate()
$ws.Range("J12").Select()
try {
    $excel.ActiveWindow.ScrollColumn = 4
    $excel.ActiveWindow.ScrollRow = 1
} catch {
}
